$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 253.5
$ws.Range("I9").Value = 269.76923
$ws.Range("J9").Value = 183
$ws.Range("K9").Value = 269.76923
$ws.Range("L9").Value = 183
$ws.Range("M9").Value = -100.76923
$ws.Range("N9").Value = -521
$ws.Range("H11").Value = 1561.5714
$ws.Range("I11").Value = 1561.5714
$ws.Range("K11").Value = 1561.5714
$ws.Range("M11").Value = -1421.5714
$ws.Range("H15").Value = 2741.2144
$ws.Range("I15").Value = 2741.2144
$ws.Range("K15").Value = 8223.643199999999
$ws.Range("M15").Value = -8054.643199999999
$ws.Range("H55").Value = 944.5454999999999
$ws.Range("I55").Value = 849.25
$ws.Range("J55").Value = 999
$ws.Range("K55").Value = 849.25
$ws.Range("L55").Value = 999
$ws.Range("M55").Value = -635.25
$ws.Range("N55").Value = -1427
$ws.Range("H112").Value = 1579
$ws.Range("I112").Value = 1216.3334
$ws.Range("K112").Value = 3649.0002
$ws.Range("M112").Value = -2541.0002
$ws.Range("H135").Value = 1498.5
$ws.Range("I135").Value = 1498.5
$ws.Range("K135").Value = 13486.5
$ws.Range("M135").Value = -10951.5
$ws.Range("H138").Value = 2632.111
$ws.Range("I138").Value = 2167.2942
$ws.Range("K138").Value = 6501.882599999999
$ws.Range("M138").Value = -1361.882599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5669.5
$ws.Range("I61").Value = 6497.5
$ws.Range("J61").Value = 4013.5
$ws.Range("K61").Value = 6497.5
$ws.Range("L61").Value = 4013.5
$ws.Range("M61").Value = -6285.5
$ws.Range("N61").Value = -4437.5
$ws.Range("H63").Value = 8499.166999999999
$ws.Range("I63").Value = 8499.5
$ws.Range("K63").Value = 8499.5
$ws.Range("M63").Value = -7813.5
$ws.Range("H66").Value = 8499.166999999999
$ws.Range("I66").Value = 8499.5
$ws.Range("K66").Value = 42497.5
$ws.Range("M66").Value = -39065.5
$ws.Range("H88").Value = 3572.7144
$ws.Range("J88").Value = 3751.5
$ws.Range("L88").Value = 3751.5
$ws.Range("N88").Value = -4563.5
$ws.Range("H91").Value = 3572.7144
$ws.Range("J91").Value = 3751.5
$ws.Range("L91").Value = 3751.5
$ws.Range("N91").Value = -6559.5
$ws.Range("H102").Value = 2355.4
$ws.Range("I102").Value = 2142.7144
$ws.Range("K102").Value = 2142.7144
$ws.Range("M102").Value = -520.7143999999998
$ws.Range("H122").Value = 2038.2222
$ws.Range("I122").Value = 1783.6818
$ws.Range("J122").Value = 3158.2
$ws.Range("K122").Value = 5351.0454
$ws.Range("L122").Value = 9474.599999999999
$ws.Range("M122").Value = -2901.0454
$ws.Range("N122").Value = -14374.6
$ws.Range("H132").Value = 66217
$ws.Range("I132").Value = 66217
$ws.Range("K132").Value = 198651
$ws.Range("M132").Value = -196121
$ws.Range("H136").Value = 5669.5
$ws.Range("I136").Value = 6497.5
$ws.Range("J136").Value = 4013.5
$ws.Range("K136").Value = 19492.5
$ws.Range("L136").Value = 12040.5
$ws.Range("M136").Value = -16942.5
$ws.Range("N136").Value = -17140.5
$ws.Range("H138").Value = 99999
$ws.Range("J138").Value = 99999
$ws.Range("L138").Value = 99999
$ws.Range("N138").Value = -110279

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4064.1304
$ws.Range("I86").Value = 3109.1667
$ws.Range("J86").Value = 5105.909
$ws.Range("K86").Value = 3109.1667
$ws.Range("L86").Value = 5105.909
$ws.Range("M86").Value = -1986.1667
$ws.Range("N86").Value = -7351.909
$ws.Range("H89").Value = 4064.1304
$ws.Range("I89").Value = 3109.1667
$ws.Range("J89").Value = 5105.909
$ws.Range("K89").Value = 15545.8335
$ws.Range("L89").Value = 25529.545
$ws.Range("M89").Value = -9929.833500000001
$ws.Range("N89").Value = -36761.545
$ws.Range("H134").Value = 3928.5625
$ws.Range("I134").Value = 3277.75
$ws.Range("J134").Value = 5881
$ws.Range("K134").Value = 9833.25
$ws.Range("L134").Value = 17643
$ws.Range("M134").Value = -7298.25
$ws.Range("N134").Value = -22713

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 691.7143
$ws.Range("J22").Value = 854.125
$ws.Range("L22").Value = 854.125
$ws.Range("N22").Value = -1554.125
$ws.Range("H31").Value = 4891.222
$ws.Range("I31").Value = 3661.3333
$ws.Range("J31").Value = 5506.1665
$ws.Range("K31").Value = 3661.3333
$ws.Range("L31").Value = 5506.1665
$ws.Range("M31").Value = -3366.3333
$ws.Range("N31").Value = -6096.1665
$ws.Range("H34").Value = 4891.222
$ws.Range("I34").Value = 3661.3333
$ws.Range("J34").Value = 5506.1665
$ws.Range("K34").Value = 3661.3333
$ws.Range("L34").Value = 5506.1665
$ws.Range("M34").Value = -3459.3333
$ws.Range("N34").Value = -5910.1665
$ws.Range("H58").Value = 69193.47
$ws.Range("J58").Value = 3583.8333
$ws.Range("L58").Value = 3583.8333
$ws.Range("N58").Value = -3989.8333
$ws.Range("H86").Value = 6265.4165
$ws.Range("J86").Value = 5599.3335
$ws.Range("L86").Value = 5599.3335
$ws.Range("N86").Value = -7845.3335
$ws.Range("H89").Value = 6265.4165
$ws.Range("J89").Value = 5599.3335
$ws.Range("L89").Value = 27996.6675
$ws.Range("N89").Value = -39228.6675
$ws.Range("H105").Value = 1121.3
$ws.Range("I105").Value = 912.55554
$ws.Range("K105").Value = 912.55554
$ws.Range("M105").Value = 834.44446
$ws.Range("H136").Value = 69193.47
$ws.Range("J136").Value = 3583.8333
$ws.Range("L136").Value = 10751.4999
$ws.Range("N136").Value = -15851.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 69162.8
$ws.Range("I132").Value = 79360.92
$ws.Range("K132").Value = 238082.76
$ws.Range("M132").Value = -235552.76

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 54623.19
$ws.Range("I22").Value = 112009.8
$ws.Range("K22").Value = 112009.8
$ws.Range("M22").Value = -111714.8
$ws.Range("H27").Value = 54623.19
$ws.Range("I27").Value = 112009.8
$ws.Range("K27").Value = 112009.8
$ws.Range("M27").Value = -111902.8
$ws.Range("H46").Value = 9905.625
$ws.Range("J46").Value = 3356.4285
$ws.Range("L46").Value = 3356.4285
$ws.Range("N46").Value = -3732.4285
$ws.Range("H68").Value = 3526.2273
$ws.Range("I68").Value = 2093.4167
$ws.Range("J68").Value = 5245.6
$ws.Range("K68").Value = 2093.4167
$ws.Range("L68").Value = 5245.6
$ws.Range("M68").Value = -1344.4167
$ws.Range("N68").Value = -6743.6
$ws.Range("H71").Value = 3526.2273
$ws.Range("I71").Value = 2093.4167
$ws.Range("J71").Value = 5245.6
$ws.Range("K71").Value = 10467.0835
$ws.Range("L71").Value = 26228
$ws.Range("M71").Value = -6723.083500000001
$ws.Range("N71").Value = -33716
$ws.Range("H74").Value = 78000
$ws.Range("I74").Value = 78000
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 78000
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -77002
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 78000
$ws.Range("I77").Value = 78000
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 234000
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -229008
$ws.Range("N77").ClearContents()
$ws.Range("H112").Value = 50128.332
$ws.Range("J112").Value = 50128.332
$ws.Range("L112").Value = 50128.332
$ws.Range("N112").Value = -53082.332
$ws.Range("H121").Value = 3420
$ws.Range("J121").Value = 3420
$ws.Range("L121").Value = 3420
$ws.Range("N121").Value = -6914

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 761499.25
$ws.Range("J29").Value = 15332.333
$ws.Range("L29").Value = 15332.333
$ws.Range("N29").Value = -15912.333
$ws.Range("H107").Value = 1700.0625
$ws.Range("I107").Value = 1207.8889
$ws.Range("J107").Value = 2332.8572
$ws.Range("K107").Value = 3623.6667
$ws.Range("L107").Value = 6998.571599999999
$ws.Range("M107").Value = -1703.6667
$ws.Range("N107").Value = -10838.5716
$ws.Range("H112").Value = 25000
$ws.Range("J112").Value = 25000
$ws.Range("L112").Value = 25000
$ws.Range("N112").Value = -27954
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280
